# fall 24 week 8 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.27
$ws.Range("B3").Value = 1.57
$ws.Range("C4").Value = 1.46
$ws.Range("F4").Value = 1.08
$ws.Range("D6").Value = 1.52
$ws.Range("E6").Value = 1.33
$ws.Range("G7").Value = 1.14
